$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "SOURCE" column header in H1, matching the bold header style used by G1
$ws.Range("H1").Value = "SOURCE"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate SOURCE = 1 for every data row (2-12)
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}

# Resize columns to fit the new layout
# (Excel stores column widths snapped to pixel granularity; the values below
# are chosen so the saved width matches the target as closely as possible.)
$ws.Columns.Item(1).ColumnWidth = 13.33
$ws.Columns.Item(2).ColumnWidth = 12.33
$ws.Columns.Item(6).ColumnWidth = 22.17
$ws.Columns.Item(7).ColumnWidth = 23.17

# Header row height
$ws.Rows.Item(1).RowHeight = 14

# Update selection to reflect where the editor ended up
$ws.Range("I6").Select() | Out-Null
